$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 11455.3125
$ws.Range("I76").Value = 29995.25
$ws.Range("J76").Value = 5275.3335
$ws.Range("K76").Value = 29995.25
$ws.Range("L76").Value = 5275.3335
$ws.Range("M76").Value = -29680.25
$ws.Range("N76").Value = -5905.3335

$ws.Range("H79").Value = 11455.3125
$ws.Range("I79").Value = 29995.25
$ws.Range("J79").Value = 5275.3335
$ws.Range("K79").Value = 29995.25
$ws.Range("L79").Value = 5275.3335
$ws.Range("M79").Value = -28903.25
$ws.Range("N79").Value = -7459.3335

$ws.Range("H141").Value = 3693.4075
$ws.Range("I141").Value = 1009.85
$ws.Range("J141").Value = 11360.714
$ws.Range("K141").Value = 3029.55
$ws.Range("L141").Value = 34082.142
$ws.Range("M141").Value = 2150.45
$ws.Range("N141").Value = -44442.142

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 320311.84
$ws.Range("I61").Value = 235444.58
$ws.Range("J61").Value = 530246.6
$ws.Range("K61").Value = 235444.58
$ws.Range("L61").Value = 530246.6
$ws.Range("M61").Value = -235232.58
$ws.Range("N61").Value = -530670.6

$ws.Range("H63").Value = 3188
$ws.Range("I63").Value = 3118.9092
$ws.Range("J63").Value = 3340
$ws.Range("K63").Value = 3118.9092
$ws.Range("L63").Value = 3340
$ws.Range("M63").Value = -2432.9092
$ws.Range("N63").Value = -4712

$ws.Range("H66").Value = 3188
$ws.Range("I66").Value = 3118.9092
$ws.Range("J66").Value = 3340
$ws.Range("K66").Value = 15594.546
$ws.Range("L66").Value = 16700
$ws.Range("M66").Value = -12162.546
$ws.Range("N66").Value = -23564

$ws.Range("H132").Value = 18825.242
$ws.Range("I132").Value = 24460.2
$ws.Range("J132").Value = 3909.1765
$ws.Range("K132").Value = 73380.60000000001
$ws.Range("L132").Value = 11727.5295
$ws.Range("M132").Value = -70850.60000000001
$ws.Range("N132").Value = -16787.5295

$ws.Range("H136").Value = 320311.84
$ws.Range("I136").Value = 235444.58
$ws.Range("J136").Value = 530246.6
$ws.Range("K136").Value = 706333.74
$ws.Range("L136").Value = 1590739.8
$ws.Range("M136").Value = -703783.74
$ws.Range("N136").Value = -1595839.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1674.7894
$ws.Range("I105").Value = 1594
$ws.Range("J105").Value = 1977.75
$ws.Range("K105").Value = 1594
$ws.Range("L105").Value = 1977.75
$ws.Range("M105").Value = 153
$ws.Range("N105").Value = -5471.75

$ws.Range("H134").Value = 2565.1128
$ws.Range("I134").Value = 2139.9387
$ws.Range("J134").Value = 4167.6924
$ws.Range("K134").Value = 6419.8161
$ws.Range("L134").Value = 12503.0772
$ws.Range("M134").Value = -3884.8161
$ws.Range("N134").Value = -17573.0772

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1000
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 1000
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 1000
$ws.Range("M16").ClearContents()
$ws.Range("N16").Value = -1574

$ws.Range("H58").Value = 8513.733
$ws.Range("I58").Value = 10930.6
$ws.Range("J58").Value = 3680
$ws.Range("K58").Value = 10930.6
$ws.Range("L58").Value = 3680
$ws.Range("M58").Value = -10727.6
$ws.Range("N58").Value = -4086

$ws.Range("H113").Value = 1000
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 1000
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 1000
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = -5340

$ws.Range("H132").Value = 3269.8333
$ws.Range("I132").Value = 1412.4445
$ws.Range("J132").Value = 5127.222
$ws.Range("K132").Value = 4237.333500000001
$ws.Range("L132").Value = 15381.666
$ws.Range("M132").Value = -1707.333500000001
$ws.Range("N132").Value = -20441.666

$ws.Range("H134").Value = 1780.0588
$ws.Range("I134").Value = 1141.2174
$ws.Range("J134").Value = 3115.818
$ws.Range("K134").Value = 3423.6522
$ws.Range("L134").Value = 9347.454000000002
$ws.Range("M134").Value = -888.6522
$ws.Range("N134").Value = -14417.454

$ws.Range("H136").Value = 8513.733
$ws.Range("I136").Value = 10930.6
$ws.Range("J136").Value = 3680
$ws.Range("K136").Value = 32791.8
$ws.Range("L136").Value = 11040
$ws.Range("M136").Value = -30241.8
$ws.Range("N136").Value = -16140

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5799.4585
$ws.Range("I70").Value = 4523.375
$ws.Range("J70").Value = 6437.5
$ws.Range("K70").Value = 4523.375
$ws.Range("L70").Value = 6437.5
$ws.Range("M70").Value = -4253.375
$ws.Range("N70").Value = -6977.5

$ws.Range("H73").Value = 5799.4585
$ws.Range("I73").Value = 4523.375
$ws.Range("J73").Value = 6437.5
$ws.Range("K73").Value = 4523.375
$ws.Range("L73").Value = 6437.5
$ws.Range("M73").Value = -3587.375
$ws.Range("N73").Value = -8309.5

$ws.Range("H80").Value = 6177.4443
$ws.Range("I80").Value = 8572.182000000001
$ws.Range("K80").Value = 8572.182000000001
$ws.Range("M80").Value = -7574.182000000001

$ws.Range("H83").Value = 6177.4443
$ws.Range("I83").Value = 8572.182000000001
$ws.Range("K83").Value = 42860.91
$ws.Range("M83").Value = -37868.91

$ws.Range("H132").Value = 3797.4167
$ws.Range("I132").Value = 3794.2258
$ws.Range("J132").Value = 3803.2354
$ws.Range("K132").Value = 11382.6774
$ws.Range("L132").Value = 11409.7062
$ws.Range("M132").Value = -8852.6774
$ws.Range("N132").Value = -16469.7062

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2784.122
$ws.Range("I132").Value = 1899.3549
$ws.Range("J132").Value = 5526.9
$ws.Range("K132").Value = 5698.0647
$ws.Range("L132").Value = 16580.7
$ws.Range("M132").Value = -3168.0647
$ws.Range("N132").Value = -21640.7

$ws.Range("H136").Value = 3104.647
$ws.Range("I136").Value = 2125.9443
$ws.Range("J136").Value = 5453.533
$ws.Range("K136").Value = 6377.8329
$ws.Range("L136").Value = 16360.599
$ws.Range("M136").Value = -3827.8329
$ws.Range("N136").Value = -21460.599

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1672.125
$ws.Range("I132").Value = 969.7406999999999
$ws.Range("J132").Value = 3130.923
$ws.Range("K132").Value = 2909.2221
$ws.Range("L132").Value = 9392.769
$ws.Range("M132").Value = -379.2221
$ws.Range("N132").Value = -14452.769
